# Update "想去人数" (people interested) counts that changed between scrapes.
# 展览  (sheet1): F3 3087 -> 3088 ; F5 101 -> 102
# 演出  (sheet2): F2 119 -> 120
# 全部类型 (sheet4): F3 119 -> 120 ; F7 3087 -> 3088 ; F10 101 -> 102

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 3088
$wsExhibition.Range("F5").Value = 102

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 120

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 120
$wsAll.Range("F7").Value = 3088
$wsAll.Range("F10").Value = 102
